$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '29.105.15'
$ws.Range("E2").NumberFormat = '@'
$ws.Range("E2").Value = '  -2.84%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.846.44'
$ws.Range("E3").NumberFormat = '@'
$ws.Range("E3").Value = '  -1.82%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").NumberFormat = '@'
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '0.7026'
$ws.Range("E5").NumberFormat = '@'
$ws.Range("E5").Value = '  -5.11%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '237.02'
$ws.Range("E6").NumberFormat = '@'
$ws.Range("E6").Value = '  -2.41%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").NumberFormat = '@'
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.3023'
$ws.Range("E8").NumberFormat = '@'
$ws.Range("E8").Value = '  -4.34%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.07392'
$ws.Range("E9").NumberFormat = '@'
$ws.Range("E9").Value = '  +2.48%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '23.26'
$ws.Range("E10").NumberFormat = '@'
$ws.Range("E10").Value = '  -6.15%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.08114'
$ws.Range("E11").NumberFormat = '@'
$ws.Range("E11").Value = '  -2.80%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.868.55'
$ws.Range("E12").NumberFormat = '@'
$ws.Range("E12").Value = '  -8.70%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.7232'
$ws.Range("E13").NumberFormat = '@'
$ws.Range("E13").Value = '  -4.02%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '5.202'
$ws.Range("E14").NumberFormat = '@'
$ws.Range("E14").Value = '  -3.65%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '88.95'
$ws.Range("E15").NumberFormat = '@'
$ws.Range("E15").Value = '  -3.70%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '29.122.76'
$ws.Range("E16").NumberFormat = '@'
$ws.Range("E16").Value = '  -2.79%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '5.766'
$ws.Range("E17").NumberFormat = '@'
$ws.Range("E17").Value = '  -6.28%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '241.04'
$ws.Range("E18").NumberFormat = '@'
$ws.Range("E18").Value = '  -3.27%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '0.000007645'
$ws.Range("E19").NumberFormat = '@'
$ws.Range("E19").Value = '  -2.72%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '12.99'
$ws.Range("E20").NumberFormat = '@'
$ws.Range("E20").Value = '  -4.39%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").NumberFormat = '@'
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '2.090.92'
$ws.Range("E22").NumberFormat = '@'
$ws.Range("E22").Value = '  -2.50%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '0.9994'
$ws.Range("E23").NumberFormat = '@'
$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '7.573'
$ws.Range("E24").NumberFormat = '@'
$ws.Range("E24").Value = '  -5.35%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '0.1474'
$ws.Range("E25").NumberFormat = '@'
$ws.Range("E25").Value = '  -5.55%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '161.68'
$ws.Range("E26").NumberFormat = '@'
$ws.Range("E26").Value = '  -2.68%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '8.957'
$ws.Range("E27").NumberFormat = '@'
$ws.Range("E27").Value = '  -3.78%  '

$ws.Range("E28").NumberFormat = '@'
$ws.Range("E28").Value = '  -3.70%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '1.929'

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '1.373'
$ws.Range("E30").NumberFormat = '@'
$ws.Range("E30").Value = '  -8.36%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '4.443'
$ws.Range("E31").NumberFormat = '@'
$ws.Range("E31").Value = '  -3.43%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.487'

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '3.998'
$ws.Range("E33").NumberFormat = '@'
$ws.Range("E33").Value = '  -5.20%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.05198'
$ws.Range("E34").NumberFormat = '@'
$ws.Range("E34").Value = '  -3.26%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.182'
$ws.Range("E35").NumberFormat = '@'
$ws.Range("E35").Value = '  -5.28%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.7085'
$ws.Range("E36").NumberFormat = '@'
$ws.Range("E36").Value = '  -6.21%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.9999'
$ws.Range("E37").NumberFormat = '@'
$ws.Range("E37").Value = '  -0.88%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '2.647'
$ws.Range("E38").NumberFormat = '@'
$ws.Range("E38").Value = '  -2.34%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.01867'
$ws.Range("E39").NumberFormat = '@'
$ws.Range("E39").Value = '  -5.06%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '2.672'
$ws.Range("E40").NumberFormat = '@'
$ws.Range("E40").Value = '  -3.17%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '0.9017'
$ws.Range("E41").NumberFormat = '@'
$ws.Range("E41").Value = '  +4.81%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.4276'
$ws.Range("E42").NumberFormat = '@'
$ws.Range("E42").Value = '  -5.99%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '5.881'
$ws.Range("E43").NumberFormat = '@'
$ws.Range("E43").Value = '  -4.72%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '70.02'
$ws.Range("E44").NumberFormat = '@'
$ws.Range("E44").Value = '  -4.00%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '1.048.07'
$ws.Range("E45").NumberFormat = '@'
$ws.Range("E45").Value = '  -6.39%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.9998'
$ws.Range("E46").NumberFormat = '@'
$ws.Range("E46").Value = '  -0.20%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '101.49'
$ws.Range("E47").NumberFormat = '@'
$ws.Range("E47").Value = '  -3.43%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '1.748'
$ws.Range("E48").NumberFormat = '@'
$ws.Range("E48").Value = '  -6.48%  '

$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '7.104'
$ws.Range("E49").NumberFormat = '@'
$ws.Range("E49").Value = '  -6.89%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '1.984.04'
$ws.Range("E50").NumberFormat = '@'
$ws.Range("E50").Value = '  -3.54%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '9.176'
$ws.Range("E51").NumberFormat = '@'
$ws.Range("E51").Value = '  -3.91%  '

